$wb = $excel.ActiveWorkbook

# --- Sheet "pacmed_table" (sheet1): selection only changes (no data edits) ---
$ws1 = $wb.Worksheets.Item("pacmed_table")
$ws1.Range("B20").Select()

# --- Sheet "search_strings" (sheet2): add the "Time Series" column ---
$ws2 = $wb.Worksheets.Item("search_strings")
$ws2.Activate()

# Drop the huge block of trailing empty formatted rows before touching columns
$ws2.Range("A1048513:A1048576").EntireRow.Delete()

# Insert a new column before the existing "Pacmed ontology" column (F),
# shifting Pacmed ontology / Snowmed CT / Notes one column to the right.
$ws2.Columns.Item(6).Insert()

# Header text for the new column (inherits the bold header style already in
# place on row 1 from the column insert, matching its neighbours)
$ws2.Range("F1").Value2 = "Time Series"

# New boolean "TRUE" values for data rows 2 & 3, with the custom boolean format
$ws2.Range("F2").Value2 = $true
$ws2.Range("F3").Value2 = $true
$ws2.Range("F2:F3").NumberFormat = '"TRUE";"TRUE";"FALSE"'

# Row heights for rows 1 & 2 now match row 3's existing 13.8 height
$ws2.Rows.Item(1).RowHeight = 13.8
$ws2.Rows.Item(2).RowHeight = 13.8

# Give the new column roughly the same width as its former position
$ws2.Columns.Item(6).ColumnWidth = 46.25

# Final selection/view state on the sheet
$ws2.Range("F4").Select()

Write-Output "time series column added"
